$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 417.6
$ws.Range("I11").Value = 417.6
$ws.Range("K11").Value = 417.6
$ws.Range("M11").Value = -277.6
$ws.Range("H121").Value = 6142.6665
$ws.Range("J121").Value = 6142.6665
$ws.Range("L121").Value = 18427.9995
$ws.Range("N121").Value = -21921.9995
$ws.Range("H138").Value = 5754.6
$ws.Range("J138").Value = 5227.871
$ws.Range("L138").Value = 15683.613
$ws.Range("N138").Value = -25963.613

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18884.707
$ws.Range("I32").Value = 8581
$ws.Range("K32").Value = 8581
$ws.Range("M32").Value = -8294
$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2000
$ws.Range("K33").Value = 2000
$ws.Range("M33").Value = -1671
$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 60000
$ws.Range("L44").Value = 60000
$ws.Range("N44").Value = -60976
$ws.Range("H45").Value = 2552.875
$ws.Range("I45").Value = 1605.75
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 1605.75
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -1228.75
$ws.Range("N45").Value = -4254
$ws.Range("H80").Value = 93258.336
$ws.Range("J80").Value = 93258.336
$ws.Range("L80").Value = 93258.336
$ws.Range("N80").Value = -95254.336
$ws.Range("H83").Value = 93258.336
$ws.Range("J83").Value = 93258.336
$ws.Range("L83").Value = 279775.008
$ws.Range("N83").Value = -289759.008
$ws.Range("H122").Value = 718535.9
$ws.Range("I122").Value = 1670250.4
$ws.Range("K122").Value = 5010751.199999999
$ws.Range("M122").Value = -5008301.199999999
$ws.Range("H134").Value = 97150
$ws.Range("I134").Value = 97000
$ws.Range("J134").Value = 97200
$ws.Range("K134").Value = 97000
$ws.Range("L134").Value = 97200
$ws.Range("M134").Value = -91930
$ws.Range("N134").Value = -107340

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2873.9333
$ws.Range("I134").Value = 1275.1
$ws.Range("K134").Value = 3825.3
$ws.Range("M134").Value = -1290.3

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4626.7383
$ws.Range("I31").Value = 4032.1
$ws.Range("K31").Value = 4032.1
$ws.Range("M31").Value = -3737.1
$ws.Range("H34").Value = 4626.7383
$ws.Range("I34").Value = 4032.1
$ws.Range("K34").Value = 4032.1
$ws.Range("M34").Value = -3830.1
$ws.Range("H94").Value = 1086
$ws.Range("I94").Value = 680
$ws.Range("J94").Value = 1136.75
$ws.Range("K94").Value = 680
$ws.Range("L94").Value = 1136.75
$ws.Range("M94").Value = -229
$ws.Range("N94").Value = -2038.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 246
$ws.Range("I2").Value = 214.42857
$ws.Range("J2").Value = 356.5
$ws.Range("K2").Value = 1286.57142
$ws.Range("L2").Value = 2139
$ws.Range("M2").Value = -1173.57142
$ws.Range("N2").Value = -2365
$ws.Range("H51").Value = 396
$ws.Range("I51").Value = 396
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1188
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -728
$ws.Range("N51").ClearContents()
$ws.Range("H107").Value = 1655.4166
$ws.Range("J107").Value = 1163.1428
$ws.Range("L107").Value = 3489.4284
$ws.Range("N107").Value = -7329.428400000001
$ws.Range("H113").Value = 3914.6667
$ws.Range("J113").Value = 4147.5
$ws.Range("L113").Value = 12442.5
$ws.Range("N113").Value = -16782.5
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H131").Value = 1455.6666
$ws.Range("I131").Value = 979.1667
$ws.Range("J131").Value = 1591.8096
$ws.Range("K131").Value = 2937.5001
$ws.Range("L131").Value = 4775.4288
$ws.Range("M131").Value = 2102.4999
$ws.Range("N131").Value = -14855.4288

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4993.3335
$ws.Range("I126").Value = 4990.4
$ws.Range("J126").Value = 4997
$ws.Range("K126").Value = 14971.2
$ws.Range("L126").Value = 14991
$ws.Range("M126").Value = -12501.2
$ws.Range("N126").Value = -19931
$ws.Range("H132").Value = 3608.3333
$ws.Range("I132").Value = 2375.2
$ws.Range("K132").Value = 7125.599999999999
$ws.Range("M132").Value = -4595.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5052499.5
$ws.Range("I2").Value = 10005000
$ws.Range("J2").Value = 99998.5
$ws.Range("K2").Value = 10005000
$ws.Range("L2").Value = 99998.5
$ws.Range("M2").Value = -10004888
$ws.Range("N2").Value = -100222.5
$ws.Range("H22").Value = 3133
$ws.Range("J22").Value = 4449.5
$ws.Range("L22").Value = 4449.5
$ws.Range("N22").Value = -5039.5
$ws.Range("H27").Value = 3133
$ws.Range("J27").Value = 4449.5
$ws.Range("L27").Value = 4449.5
$ws.Range("N27").Value = -4663.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 61055
$ws.Range("J97").Value = 61055
$ws.Range("L97").Value = 61055
$ws.Range("N97").Value = -63037
$ws.Range("H98").Value = 28350
$ws.Range("J98").Value = 28350
$ws.Range("L98").Value = 28350
$ws.Range("N98").Value = -34340
$ws.Range("H99").Value = 13172.333
$ws.Range("I99").Value = 13172.333
$ws.Range("K99").Value = 13172.333
$ws.Range("M99").Value = -10177.333
$ws.Range("H132").Value = 5090.1
$ws.Range("I132").Value = 3700.3333
$ws.Range("J132").Value = 6227.1816
$ws.Range("K132").Value = 11100.9999
$ws.Range("L132").Value = 18681.5448
$ws.Range("M132").Value = -8570.999899999999
$ws.Range("N132").Value = -23741.5448

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 73683.07000000001
$ws.Range("I14").Value = 92782.55
$ws.Range("J14").Value = 3651.6667
$ws.Range("K14").Value = 92782.55
$ws.Range("L14").Value = 3651.6667
$ws.Range("M14").Value = -92614.55
$ws.Range("N14").Value = -3987.6667
$ws.Range("H25").Value = 37916
$ws.Range("J25").Value = 37916
$ws.Range("L25").Value = 37916
$ws.Range("N25").Value = -38502
